$d = $word.ActiveDocument

# Helper: replace the trailing underline-run text of a paragraph.
# $paraIndex  - 1-based paragraph index
# $labelLen   - length (in characters) of the bold label run(s) that precede
#               the underline run inside the paragraph
# $newText    - full new text for the underline run (replaces whatever is
#               currently there, including any existing spaces/placeholder)
function Set-UnderlineSuffix($paraIndex, $labelLen, $newText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $full = $p.Range.Duplicate
    $full.MoveEnd(1, -1)          # exclude the trailing paragraph mark
    $start = $full.Start + $labelLen
    $end = $full.End
    $r = $d.Range($start, $end)
    $r.Text = $newText
}

# 1. ΕΠΩΝΥΜΟ:  " " -> " ασδασφλδφηδκσξηω"
Set-UnderlineSuffix 2 9 " ασδασφλδφηδκσξηω"

# 2. ΟΝΟΜΑ:    "" -> "ωδσφσδγδσαφγ"
Set-UnderlineSuffix 3 10 "ωδσφσδγδσαφγ"

# 3. ΟΝΟΜΑ ΠΑΤΕΡΑ:    "" -> "γφδφγδφγσδγφδσ"
Set-UnderlineSuffix 4 17 "γφδφγδφγσδγφδσ"

# 4. ΟΝΟΜΑ ΜΗΤΕΡΑΣ:    " " -> " φγδγδφγφδδσ"
Set-UnderlineSuffix 5 18 " φγδγδφγφδδσ"

# 5. ΗΜΕΡΟΜΗΝΙΑ ΓΕΝΝΗΣΗΣ:    "  " -> "  φδγσγφγδφδσγ"
Set-UnderlineSuffix 8 24 "  φδγσγφγδφδσγ"

# 6.  ΤΟΠΟΘΕΣΙΑ ΓΕΝΝΗΣΗΣ:    "  " -> "  γφδσφγφδγδσγδφσ"
Set-UnderlineSuffix 9 24 "  γφδσφγφδγδσγδφσ"

# 7.  ΤΟΠΟΣ ΚΑΤΟΙΚΙΑΣ:    "  " -> "  φγδδσφγφδφσγφδσ"
Set-UnderlineSuffix 10 21 "  φγδδσφγφδφσγφδσ"

# 8.  ΑΡΙΘΜΟΣ Δ.Α.Τ:    "  " -> "  γφδσγδφσ"
Set-UnderlineSuffix 11 19 "  γφδσγδφσ"

# 9.  ΕΚΔΟΘΕΝ:    "  " -> "  γφδσσφδγδσφ"
Set-UnderlineSuffix 12 13 "  γφδσσφδγδσφ"

# 10.  ΑΠΟ:    " " -> " γφδσφγδφσ"
Set-UnderlineSuffix 13 9 " γφδσφγδφσ"

# 11.  Α.Φ.Μ:    "   από Δ.Ο.Υ   " -> "  γφσδγδφσγ από Δ.Ο.Υ γφδδφσγδφγδσ  "
Set-UnderlineSuffix 14 11 "  γφσδγδφσγ από Δ.Ο.Υ γφδδφσγδφγδσ  "

# 12. Date: "Θέρμη, 04/10/2025" -> "Θέρμη, 09/11/2025"
$d.Content.Find.Execute("Θέρμη, 04/10/2025", $true, $true, $false, $false, $false, $true, 1, $false, "Θέρμη, 09/11/2025", 2) | Out-Null

# 13. Signature run (last paragraph, trailing " " before the final <w:br/>):
#     " " -> " Υ/Α Αθανασιάδης Γρηγορίου"
$p15 = $d.Paragraphs.Item($d.Paragraphs.Count)
$p15full = $p15.Range.Duplicate
$p15full.MoveEnd(1, -1)
$sigR = $p15full.Duplicate
$sigR.MoveEnd(1, -1)          # drop the trailing <w:br/> char too
$sigR.Collapse(0)             # wdCollapseEnd -> sits right before the <w:br/>
$sigR.MoveStart(1, -1)        # extend back to cover the single-space run
$sigR.Text = " Υ/Α Αθανασιάδης Γρηγορίου"
